$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = -0.03986657931951261
$ws.Range("C2").Value = 1.166715395314321
$ws.Range("D2").Value = 9.175165609364871
$ws.Range("E2").Value = 3.029053583112202
$ws.Range("F2").Value = 3.058341104199183

# Row 3 (Q1)
$ws.Range("B3").Value = 0.3935572759351038
$ws.Range("C3").Value = 0.9729011477362982
$ws.Range("D3").Value = 6.967390270327681
$ws.Range("E3").Value = 2.639581457414732
$ws.Range("F3").Value = 2.636048747596285

# Row 4 (Q2)
$ws.Range("B4").Value = 0.5845645237101601
$ws.Range("C4").Value = 0.9612686620598623
$ws.Range("D4").Value = 4.440372314603645
$ws.Range("E4").Value = 2.107219095064309
$ws.Range("F4").Value = 2.045067892746351

# Row 5 (Q3)
$ws.Range("B5").Value = 0.4605577599719269
$ws.Range("C5").Value = 1.029373326337189
$ws.Range("D5").Value = 4.397813778901289
$ws.Range("E5").Value = 2.097096511584836
$ws.Range("F5").Value = 2.071314277662061
$ws.Range("G5").Value = 41

# Row 6 (Q4)
$ws.Range("B6").Value = 0.5450649147966192
$ws.Range("C6").Value = 1.227243952030447
$ws.Range("D6").Value = 6.248613964862616
$ws.Range("E6").Value = 2.499722777602072
$ws.Range("F6").Value = 2.479899623702079
$ws.Range("G6").Value = 31

# Row 7 (Q5)
$ws.Range("B7").Value = 0.4748548997082894
$ws.Range("C7").Value = 1.211225564801436
$ws.Range("D7").Value = 6.526816765525255
$ws.Range("E7").Value = 2.554763543955733
$ws.Range("F7").Value = 2.554677489448803
$ws.Range("G7").Value = 29

# Row 8 (Q6)
$ws.Range("B8").Value = 0.411137621314042
$ws.Range("C8").Value = 1.221296936513468
$ws.Range("D8").Value = 6.896165483727119
$ws.Range("E8").Value = 2.626055118181474
$ws.Range("F8").Value = 2.643079106050209
$ws.Range("G8").Value = 27

# Row 9 (Q7)
$ws.Range("B9").Value = 0.4034568762020432
$ws.Range("C9").Value = 1.549620193073209
$ws.Range("D9").Value = 9.580566432589528
$ws.Range("E9").Value = 3.095249009787343
$ws.Range("F9").Value = 3.15293505810293
$ws.Range("G9").Value = 19

# Row 10 (Q8)
$ws.Range("B10").Value = 0.159026239978337
$ws.Range("C10").Value = 1.962138362601939
$ws.Range("D10").Value = 14.29507177400229
$ws.Range("E10").Value = 3.780882406793723
$ws.Range("F10").Value = 3.945508253329552
$ws.Range("G10").Value = 12

# Row 11 (Q9)
$ws.Range("B11").Value = -1.776008350391145
$ws.Range("C11").Value = 2.037774147513916
$ws.Range("D11").Value = 13.78833596834753
$ws.Range("E11").Value = 3.713264866441327
$ws.Range("F11").Value = 3.645910432883749
